$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2643290936523215
$ws.Range("C2").Value = 0.0585737639768098
$ws.Range("E2").Value = 0.1749528900206911
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.3168819009209329
$ws.Range("H2").Value = 0.5047245492601675
$ws.Range("I2").Value = 0.3922376018453271
$ws.Range("K2").Value = 0.2801417449290113
$ws.Range("M2").Value = 0.2115747923263669
$ws.Range("O2").Value = 1.573953420813424
$ws.Range("B3").Value = 0.2306756071028531
$ws.Range("C3").Value = 0.05480638212192446
$ws.Range("E3").Value = 0.164369629743625
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.3218032406899951
$ws.Range("H3").Value = 0.5104184799795739
$ws.Range("I3").Value = 0.3986876248649338
$ws.Range("K3").Value = 0.2455103952596005
$ws.Range("M3").Value = 0.1885901188429813
$ws.Range("O3").Value = 1.596099716611093
$ws.Range("B4").Value = 0.2099324705357617
$ws.Range("C4").Value = 0.05248045025592774
$ws.Range("E4").Value = 0.1580070184653621
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.3251053397435051
$ws.Range("H4").Value = 0.5141542175905975
$ws.Range("I4").Value = 0.4029023134959324
$ws.Range("K4").Value = 0.2241445995302058
$ws.Range("M4").Value = 0.1745133791047309
$ws.Range("O4").Value = 1.610786782893364
$ws.Range("B5").Value = 0.2014600260715156
$ws.Range("C5").Value = 0.05152947622295301
$ws.Range("E5").Value = 0.1554480970526342
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.3265213432891585
$ws.Range("H5").Value = 0.5157368275565588
$ws.Range("I5").Value = 0.4046837377835857
$ws.Range("K5").Value = 0.2154128216540272
$ws.Range("M5").Value = 0.1687861294694457
$ws.Range("O5").Value = 1.617045488359814
$ws.Range("B6").Value = 0.2000520264630268
$ws.Range("C6").Value = 0.05137138015740561
$ws.Range("E6").Value = 0.1550252316321945
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.32676071599154
$ws.Range("H6").Value = 0.516003258891967
$ws.Range("I6").Value = 0.404983400030412
$ws.Range("K6").Value = 0.2139614220102715
$ws.Range("M6").Value = 0.1678356787429323
$ws.Range("O6").Value = 1.61810125890571
$ws.Range("B7").Value = 0.2098182861229816
$ws.Range("C7").Value = 0.05246763770493601
$ws.Range("E7").Value = 0.1579723709728142
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.3251241517059249
$ws.Range("H7").Value = 0.514175317195722
$ws.Range("I7").Value = 0.40292607973117
$ws.Range("K7").Value = 0.2240269403091162
$ws.Range("M7").Value = 0.1744361022907484
$ws.Range("O7").Value = 1.610870082518218
$ws.Range("B8").Value = 0.2527422798815451
$ws.Range("C8").Value = 0.05727745458470679
$ws.Range("E8").Value = 0.1712755290796579
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.3185205079835178
$ws.Range("H8").Value = 0.5066380835681947
$ws.Range("I8").Value = 0.3944087482894574
$ws.Range("K8").Value = 0.2682223808500623
$ws.Range("M8").Value = 0.2036422686686734
$ws.Range("O8").Value = 1.581363250051858
$ws.Range("B9").Value = 0.3362604749288209
$ws.Range("C9").Value = 0.06660587944008967
$ws.Range("E9").Value = 0.1984486563631336
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.3078013198785357
$ws.Range("H9").Value = 0.4937586246828687
$ws.Range("I9").Value = 0.3797266186880748
$ws.Range("K9").Value = 0.3540574543959849
$ws.Range("M9").Value = 0.2611992849678728
$ws.Range("O9").Value = 1.532154169864469
$ws.Range("B10").Value = 0.397196829831131
$ws.Range("C10").Value = 0.07339360534443529
$ws.Range("E10").Value = 0.219090740279313
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.3012937711673658
$ws.Range("H10").Value = 0.4854542923017675
$ws.Range("I10").Value = 0.3701740173393162
$ws.Range("K10").Value = 0.4165887087972067
$ws.Range("M10").Value = 0.3036619270975933
$ws.Range("O10").Value = 1.501291653654576
$ws.Range("B11").Value = 0.4248213126291773
$ws.Range("C11").Value = 0.07646668528660427
$ws.Range("E11").Value = 0.2286321592796909
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.2986320521281485
$ws.Range("H11").Value = 0.4819277802782125
$ws.Range("I11").Value = 0.3660968920462917
$ws.Range("K11").Value = 0.4449155111551875
$ws.Range("M11").Value = 0.323018181036467
$ws.Range("O11").Value = 1.488403868148055
$ws.Range("B12").Value = 0.4352676766906711
$ws.Range("C12").Value = 0.07762820978884122
$ws.Range("E12").Value = 0.232267233118975
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.2976672094204105
$ws.Range("H12").Value = 0.4806284961997633
$ws.Range("I12").Value = 0.3645916440329806
$ws.Range("K12").Value = 0.455624490508626
$ws.Range("M12").Value = 0.3303535754939944
$ws.Range("O12").Value = 1.48368951637697
$ws.Range("B13").Value = 0.4330185173948848
$ws.Range("C13").Value = 0.07737815298469286
$ws.Range("E13").Value = 0.2314833757394084
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.297873086556379
$ws.Range("H13").Value = 0.4809067132195253
$ws.Range("I13").Value = 0.3649141050341882
$ws.Range("K13").Value = 0.4533189207710393
$ws.Range("M13").Value = 0.3287735198494843
$ws.Range("O13").Value = 1.48469744904159
$ws.Range("B14").Value = 0.4256810341974813
$ws.Range("C14").Value = 0.07656228879024241
$ws.Range("E14").Value = 0.2289307780357959
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.2985518094169848
$ws.Range("H14").Value = 0.4818201632398171
$ws.Range("I14").Value = 0.3659722789351072
$ws.Range("K14").Value = 0.4457969051404405
$ws.Range("M14").Value = 0.3236215567853051
$ws.Range("O14").Value = 1.488012687213413
$ws.Range("B15").Value = 0.4211847154823545
$ws.Range("C15").Value = 0.07606226184877585
$ws.Range("E15").Value = 0.2273701020768613
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.2989731632592765
$ws.Range("H15").Value = 0.4823843832224526
$ws.Range("I15").Value = 0.3666254793166086
$ws.Range("K15").Value = 0.4411871202756856
$ws.Range("M15").Value = 0.3204665565227955
$ws.Range("O15").Value = 1.490064991426436
$ws.Range("B16").Value = 0.3953895215969681
$ws.Range("C16").Value = 0.07319247059278666
$ws.Range("E16").Value = 0.2184702462701722
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.3014737425430098
$ws.Range("H16").Value = 0.4856898133255072
$ws.Range("I16").Value = 0.3704458744041599
$ws.Range("K16").Value = 0.414735035267114
$ws.Range("M16").Value = 0.3023977407207923
$ws.Range("O16").Value = 1.502157106906097
$ws.Range("B17").Value = 0.3795400086528105
$ws.Range("C17").Value = 0.07142813286029082
$ws.Range("E17").Value = 0.2130493485908289
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.3030843551572886
$ws.Range("H17").Value = 0.4877819265258125
$ws.Range("I17").Value = 0.3728583604696762
$ws.Range("K17").Value = 0.3984765962592007
$ws.Range("M17").Value = 0.2913232122218687
$ws.Range("O17").Value = 1.509870466857251
$ws.Range("B18").Value = 0.3704148126850839
$ws.Range("C18").Value = 0.07041195397079036
$ws.Range("E18").Value = 0.2099456144250738
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.3040388331543795
$ws.Range("H18").Value = 0.4890088953551768
$ws.Range("I18").Value = 0.3742712173347673
$ws.Range("K18").Value = 0.3891140123917864
$ws.Range("M18").Value = 0.2849571882905693
$ws.Range("O18").Value = 1.514415372385741
$ws.Range("B19").Value = 0.367323656644345
$ws.Range("C19").Value = 0.07006765869373055
$ws.Range("E19").Value = 0.2088971800801289
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.3043668238346697
$ws.Range("H19").Value = 0.4894283858767352
$ws.Range("I19").Value = 0.3747539223064198
$ws.Range("K19").Value = 0.3859421057268833
$ws.Range("M19").Value = 0.2828024100169984
$ws.Range("O19").Value = 1.515972803653455
$ws.Range("B20").Value = 0.3812281500383108
$ws.Range("C20").Value = 0.07161609281349968
$ws.Range("E20").Value = 0.2136249389561087
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.3029099937189557
$ws.Range("H20").Value = 0.4875567707023976
$ws.Range("I20").Value = 0.3725989324852694
$ws.Range("K20").Value = 0.4002084946083926
$ws.Range("M20").Value = 0.2925017272415857
$ws.Range("O20").Value = 1.509038146509823
$ws.Range("B21").Value = 0.4278366276667214
$ws.Range("C21").Value = 0.07680198772935398
$ws.Range("E21").Value = 0.2296799406823595
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.2983512813139981
$ws.Range("H21").Value = 0.4815508801904116
$ws.Range("I21").Value = 0.3656604174962723
$ws.Range("K21").Value = 0.4480067908509398
$ws.Range("M21").Value = 0.3251346621657518
$ws.Range("O21").Value = 1.487034414114206
$ws.Range("B22").Value = 0.4582135617123413
$ws.Range("C22").Value = 0.08017850578724506
$ws.Range("E22").Value = 0.2403008194580138
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.2956231466021819
$ws.Range("H22").Value = 0.4778362793113331
$ws.Range("I22").Value = 0.3613511349084604
$ws.Range("K22").Value = 0.4791418842526696
$ws.Range("M22").Value = 0.3464947636152402
$ws.Range("O22").Value = 1.473621284200007
$ws.Range("B23").Value = 0.4420087591495587
$ws.Range("C23").Value = 0.07837758605002421
$ws.Range("E23").Value = 0.2346204784835066
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.2970561634278468
$ws.Range("H23").Value = 0.4797995583463504
$ws.Range("I23").Value = 0.3636304277486921
$ws.Range("K23").Value = 0.4625342181356018
$ws.Range("M23").Value = 0.3350915330095106
$ws.Range("O23").Value = 1.48069147312242
$ws.Range("B24").Value = 0.380464982072823
$ws.Range("C24").Value = 0.07153112183669919
$ws.Range("E24").Value = 0.2133646746605393
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.3029887337138248
$ws.Range("H24").Value = 0.4876584883158657
$ws.Range("I24").Value = 0.3727161392481939
$ws.Range("K24").Value = 0.3994255512554901
$ws.Range("M24").Value = 0.2919689178844536
$ws.Range("O24").Value = 1.509414094701967
$ws.Range("B25").Value = 0.3137393512328686
$ws.Range("C25").Value = 0.06409366667706706
$ws.Range("E25").Value = 0.1909797066214907
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.3104615160377762
$ws.Range("H25").Value = 0.4970394156165341
$ws.Range("I25").Value = 0.3834819529645177
$ws.Range("K25").Value = 0.3309284167116857
$ws.Range("M25").Value = 0.2455979453934702
$ws.Range("O25").Value = 1.54453844681214

Write-Host "applied 240 cell updates"